$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 102, shifting rows 102:235 down to 103:236
$ws.Rows.Item(102).Insert()

$ws.Range("A102").Value = 10
$ws.Range("B102").Value = "Vega Modelo de Temuco"
$ws.Range("C102").Value = "La Araucanía"
$ws.Range("D102").Value = 44546
$ws.Range("E102").Value = 9
$ws.Range("F102").Value = 100112009
$ws.Range("G102").Value = "Acelga"
$ws.Range("H102").Value = "Sin especificar"
$ws.Range("I102").Value = "Primera"
$ws.Range("J102").Value = 125
$ws.Range("K102").Value = 8000
$ws.Range("L102").Value = 8000
$ws.Range("M102").Value = 8000
$ws.Range("N102").Value = "$/docena de atados (12 kilos)"
$ws.Range("O102").Value = "Provincia de Cautín"
$ws.Range("P102").Value = 667
$ws.Range("Q102").Value = 12
$ws.Range("R102").Value = "Hortaliza"
